# Merge the "page break" paragraph with the following
# "C. Mô phỏng quá trình Truy vấn (Read Simulation) ..." paragraph into a
# single paragraph (no paragraph properties), dropping the now-redundant
# <w:lastRenderedPageBreak/> on the page-break run, per the target diff.

$d = $word.ActiveDocument

# --- locate the "C. Mô phỏng ..." paragraph robustly via Find -------------
$anchorText = "C. Mô phỏng quá trình Truy vấn"
$findRng = $d.Content
$found = $findRng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target paragraph text."
}

$paraCount = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $findRng.Start -and $p.Range.End -ge $findRng.End) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not resolve the paragraph containing the matched text."
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$breakPara = $d.Paragraphs.Item($targetIndex - 1)

# Sanity check: the preceding paragraph should just hold a manual page break.
if ($breakPara.Range.Text -notmatch [char]12) {
    throw "Preceding paragraph does not contain the expected page break."
}

$mergedRange = $d.Range($breakPara.Range.Start, $targetPara.Range.End)

# --- build the replacement OOXML for the merged paragraph -----------------
$bodyXml = '<w:p>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br w:type="page"/></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/>' +
    '<w:t>C. Mô phỏng quá trình Truy vấn (Read Simulation)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Sử dụng lệnh .explain() để chứng minh cơ chế định tuyến của MongoDB.</w:t></w:r>' +
    '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $bodyXml + '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$mergedRange.InsertXML($packageXml)

Write-Output "Merged page-break paragraph into the following paragraph."
